$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 24.916566
$ws.Range("H2").Value = 74.749698
$ws.Range("I2").Value = 0.459912889255076
$ws.Range("J2").Value = 0.459912889255076
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 23.77057966666666
$ws.Range("N2").Value = 71.31173899999999
$ws.Range("O2").Value = 0.3626243450559418
$ws.Range("P2").Value = 0.3626243450559418
$ws.Range("Q2").Value = 592.2812171227579
$ws.Range("R2").Value = 5330.530954104821
$ws.Range("S2").Value = 0.1667756102489078
$ws.Range("T2").Value = 0.1667756102489078
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 24.916566
$ws.Range("H3").Value = 74.749698
$ws.Range("I3").Value = 0.459912889255076
$ws.Range("J3").Value = 0.459912889255076
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 29.46642766666666
$ws.Range("N3").Value = 88.399283
$ws.Range("O3").Value = 0.4495155012457325
$ws.Range("P3").Value = 0.4495155012457325
$ws.Range("Q3").Value = 734.2021897407259
$ws.Range("R3").Value = 6607.819707666534
$ws.Range("S3").Value = 0.2067379729428686
$ws.Range("T3").Value = 0.2067379729428686
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 24.916566
$ws.Range("H4").Value = 74.749698
$ws.Range("I4").Value = 0.459912889255076
$ws.Range("J4").Value = 0.459912889255076
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 12.31452
$ws.Range("N4").Value = 36.94356
$ws.Range("O4").Value = 0.1878601536983258
$ws.Range("P4").Value = 0.1878601536983257
$ws.Range("Q4").Value = 306.83555033832
$ws.Range("R4").Value = 2761.519953044879
$ws.Range("S4").Value = 0.08639930606329965
$ws.Range("T4").Value = 0.08639930606329964
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 19.60300333333333
$ws.Range("H5").Value = 58.80901
$ws.Range("I5").Value = 0.3618345281251927
$ws.Range("J5").Value = 0.3618345281251927
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 23.77057966666666
$ws.Range("N5").Value = 71.31173899999999
$ws.Range("O5").Value = 0.3626243450559418
$ws.Range("P5").Value = 0.3626243450559418
$ws.Range("Q5").Value = 465.9747524409321
$ws.Range("R5").Value = 4193.772771968389
$ws.Range("S5").Value = 0.1312100087800237
$ws.Range("T5").Value = 0.1312100087800237
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 19.60300333333333
$ws.Range("H6").Value = 58.80901
$ws.Range("I6").Value = 0.3618345281251927
$ws.Range("J6").Value = 0.3618345281251927
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 29.46642766666666
$ws.Range("N6").Value = 88.399283
$ws.Range("O6").Value = 0.4495155012457325
$ws.Range("P6").Value = 0.4495155012457325
$ws.Range("Q6").Value = 577.6304797710922
$ws.Range("R6").Value = 5198.67431793983
$ws.Range("S6").Value = 0.1626502292782091
$ws.Range("T6").Value = 0.1626502292782091
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 19.60300333333333
$ws.Range("H7").Value = 58.80901
$ws.Range("I7").Value = 0.3618345281251927
$ws.Range("J7").Value = 0.3618345281251927
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 12.31452
$ws.Range("N7").Value = 36.94356
$ws.Range("O7").Value = 0.1878601536983258
$ws.Range("P7").Value = 0.1878601536983257
$ws.Range("Q7").Value = 241.4015766084
$ws.Range("R7").Value = 2172.6141894756
$ws.Range("S7").Value = 0.06797429006695987
$ws.Range("T7").Value = 0.06797429006695986
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 9.657138000000002
$ws.Range("H8").Value = 28.971414
$ws.Range("I8").Value = 0.1782525826197313
$ws.Range("J8").Value = 0.1782525826197313
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 23.77057966666666
$ws.Range("N8").Value = 71.31173899999999
$ws.Range("O8").Value = 0.3626243450559418
$ws.Range("P8").Value = 0.3626243450559418
$ws.Range("Q8").Value = 229.555768180994
$ws.Range("R8").Value = 2066.001913628946
$ws.Range("S8").Value = 0.06463872602701019
$ws.Range("T8").Value = 0.06463872602701019
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 9.657138000000002
$ws.Range("H9").Value = 28.971414
$ws.Range("I9").Value = 0.1782525826197313
$ws.Range("J9").Value = 0.1782525826197313
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 29.46642766666666
$ws.Range("N9").Value = 88.399283
$ws.Range("O9").Value = 0.4495155012457325
$ws.Range("P9").Value = 0.4495155012457325
$ws.Range("Q9").Value = 284.561358344018
$ws.Range("R9").Value = 2561.052225096162
$ws.Range("S9").Value = 0.08012729902465483
$ws.Range("T9").Value = 0.08012729902465485
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 9.657138000000002
$ws.Range("H10").Value = 28.971414
$ws.Range("I10").Value = 0.1782525826197313
$ws.Range("J10").Value = 0.1782525826197313
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 12.31452
$ws.Range("N10").Value = 36.94356
$ws.Range("O10").Value = 0.1878601536983258
$ws.Range("P10").Value = 0.1878601536983257
$ws.Range("Q10").Value = 118.92301904376
$ws.Range("R10").Value = 1070.30717139384
$ws.Range("S10").Value = 0.03348655756806623
$ws.Range("T10").Value = 0.03348655756806622
